# Update cryptos list with latest coinranking.com snapshot data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.956.38'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '1.909.30'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').Value = '''0.9987'
$ws.Range('E4').Value = '  -0.91%  '
$ws.Range('D5').Value = '''313.46'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').Value = '''0.9984'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('D7').Value = '''0.5011'
$ws.Range('E7').Value = '  +4.07%  '
$ws.Range('D8').Value = '''0.3823'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').Value = '''0.07320'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').Value = '''0.9133'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('D11').Value = '''21.21'
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '''0.07680'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '1.912.64'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '''5.515'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').Value = '''92.65'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '''0.9991'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '''0.000008760'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').Value = '''0.9978'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').Value = '27.984.86'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '''14.70'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.133.30'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '''10.86'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '''6.606'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''152.98'
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''1.846'
$ws.Range('E26').Value = '  -3.69%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '''2.214'
$ws.Range('E27').Value = '  +4.58%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''18.43'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''115.60'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '''4.924'
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '''0.09028'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = '''3.210'
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''4.853'
$ws.Range('E33').Value = '  +3.99%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''1.237'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '''0.7796'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.02088'
$ws.Range('E36').Value = '  +1.80%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''2.589'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''3.073'
$ws.Range('E38').Value = '  +2.48%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''1.095'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5569'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '''0.05289'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''6.895'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''8.537'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '''113.22'
$ws.Range('E44').Value = '  +4.34%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '''0.1523'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').Value = '''10.65'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.4840'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = '''0.9981'
$ws.Range('E48').Value = '  -0.71%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''1.641'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''67.65'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.06052'
$ws.Range('E51').Value = '  -0.51%  '
